# Export by date incident, off-bridges, maintenance
# Updates header labels, refreshes row 2 (existing incident), and appends
# four new maintenance rows (3-6) to the "Rapport maintenance" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Write a plain-text value into a cell without Excel's "looks like a
# number" auto-coercion eating leading zeros (e.g. "03250016"), and
# without leaving the cell permanently tagged with a Text number format.
function Set-TextCell($addr, $text) {
    if ($text -eq "") {
        # A bare empty string clears the cell instead of storing "", so
        # use the leading-apostrophe text marker to force an empty text
        # value to be stored.
        $ws.Range($addr).Value = "'"
    } else {
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $text
    }
    $ws.Range($addr).Style = "Normal"
}

# Write a date serial into a cell using the same numeric-date style
# already used by column B (xf index 1 / numFmtId 14), by
# copy/paste-special-formats from B2 (which already carries that style).
function Set-DateCell($addr, $serial) {
    $ws.Range($addr).Value = $serial
    $ws.Range("B2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Header relabeling (row 1)
# ---------------------------------------------------------------------
Set-TextCell "J1" "Intervenant"
Set-TextCell "K1" "Initiateur"
Set-TextCell "N1" "Cloturé par"

# ---------------------------------------------------------------------
# Row 2 : existing incident record updated (re-opened / re-dated) and a
# closing date (L2) added.
# ---------------------------------------------------------------------
Set-TextCell "A2" "03250016"
Set-TextCell "C2" "N/A"
Set-TextCell "D2" "N/A"
Set-TextCell "E2" "AIRE DE PESEES"
Set-TextCell "F2" "SALLE DE CONFERENCE ROOM"
Set-TextCell "I2" ""
Set-TextCell "K2" "0237a803-e675-49df-9d7b-25f2b329704b"
Set-TextCell "M2" "Admin User"
Set-TextCell "O2" "CLOTURE"

Set-DateCell "B2" 45734.41488488426
Set-DateCell "L2" 45734.415305752314

# ---------------------------------------------------------------------
# Row 3 : new incident
# ---------------------------------------------------------------------
Set-TextCell "A3" "03250019"
Set-TextCell "C3" "03250013"
Set-TextCell "D3" "03250013"
Set-TextCell "E3" "AIRE DE PESEES"
Set-TextCell "F3" "SALLE DE CONFERENCE ROOM"
Set-TextCell "I3" ""
Set-TextCell "J3" "HASSAN IBRAHIM  MALAM"
Set-TextCell "K3" "0237a803-e675-49df-9d7b-25f2b329704b"
Set-TextCell "M3" "Admin User"
Set-TextCell "N3" "Admin User"
Set-TextCell "O3" "CLOTURE"

Set-DateCell "B3" 45734.752714131944
Set-DateCell "L3" 45735.22049237268

# ---------------------------------------------------------------------
# Row 4 : new incident
# ---------------------------------------------------------------------
Set-TextCell "A4" "03250020"
Set-TextCell "C4" "03250012"
Set-TextCell "D4" "03250012"
Set-TextCell "E4" "AIRE DE PESEES"
Set-TextCell "F4" "SALLE DE CONFERENCE ROOM"
Set-TextCell "I4" ""
Set-TextCell "J4" "SERGE EYENGA MESSI"
Set-TextCell "K4" "0237a803-e675-49df-9d7b-25f2b329704b"
Set-TextCell "M4" "Admin User"
Set-TextCell "O4" "CLOTURE"

Set-DateCell "B4" 45734.760756006945

# ---------------------------------------------------------------------
# Row 5 : new incident
# ---------------------------------------------------------------------
Set-TextCell "A5" "03250017"
Set-TextCell "C5" "03250007"
Set-TextCell "D5" "03250007"
Set-TextCell "E5" "AIRE DE PESEES"
Set-TextCell "F5" "SALLE DE CONFERENCE ROOM"
Set-TextCell "I5" ""
Set-TextCell "J5" "Admin User"
Set-TextCell "K5" "0237a803-e675-49df-9d7b-25f2b329704b"
Set-TextCell "M5" "Admin User"
Set-TextCell "O5" "CLOTURE"

Set-DateCell "B5" 45734.414886053244

# ---------------------------------------------------------------------
# Row 6 : new incident (deleted equipment reference)
# ---------------------------------------------------------------------
Set-TextCell "A6" "03250018"
Set-TextCell "C6" "deleted__03250008__2025-03-18T14:56:20.937Z"
Set-TextCell "D6" "deleted__03250008__2025-03-18T14:56:20.937Z"
Set-TextCell "E6" "AIRE DE PESEES"
Set-TextCell "F6" "SALLE D'ATTENTE NIVEAU 1"
Set-TextCell "I6" "Something happened"
Set-TextCell "J6" "FATAHOU MOUHIDINE NCHEROLIAGNIGNI"
Set-TextCell "K6" "0237a803-e675-49df-9d7b-25f2b329704b"
Set-TextCell "M6" "Admin User"
Set-TextCell "O6" "CLOTURE"

Set-DateCell "B6" 45734.438369814816
